# Auto-update draw results: append the 2025-12-07 Pick 4 draw as a new
# row (row 82) at the bottom of the results table on the "Results" sheet,
# and grow the used range/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$newRow = 82

# Columns A ("Date") and C ("Phase") hold values that look like a date
# and a plain integer respectively ("2025-12-07" / "251207"), but in this
# workbook every column is stored as literal text. Force those two cells
# to text format before assignment so Excel keeps them as strings instead
# of silently coercing them into a date serial / number.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("C" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-12-07"
$ws.Range("B" + $newRow).Value = "Pick 4"
$ws.Range("C" + $newRow).Value = "251207"
$ws.Range("D" + $newRow).Value = "0-6-2-0"
$ws.Range("E" + $newRow).Value = "2025-12-07T21:36:59.923+04:00"
